$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# "Latest HO Xliff Generate Date" (Overview!G) / "Correspond Handoff Datetime" (de-de!H)
# shared value: 2016-10-14 08:44:53 -> 2016-10-14 08:46:50
$wsOverview.Range("G2").Value = "2016-10-14 08:46:50"
$wsOverview.Range("G3").Value = "2016-10-14 08:46:50"
$wsDeDe.Range("H2").Value = "2016-10-14 08:46:50"
$wsDeDe.Range("H3").Value = "2016-10-14 08:46:50"

# "Priority" (zh-cn!E): ht -> mt
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"

# "Correspond Handoff Datetime" (zh-cn!H): 2016-10-14 08:44:39 -> 2016-10-14 08:46:39
$wsZhCn.Range("H2").Value = "2016-10-14 08:46:39"
$wsZhCn.Range("H3").Value = "2016-10-14 08:46:39"

# "Correspond Handback DateTime" (zh-cn!K): 2016-10-14 08:45:41 -> 2016-10-14 08:47:21
$wsZhCn.Range("K2").Value = "2016-10-14 08:47:21"
$wsZhCn.Range("K3").Value = "2016-10-14 08:47:21"

# "Correspond Handback DateTime" (de-de!K): 2016-10-14 08:45:57 -> 2016-10-14 08:47:37
$wsDeDe.Range("K2").Value = "2016-10-14 08:47:37"
$wsDeDe.Range("K3").Value = "2016-10-14 08:47:37"
